# Add "Area" / "Atotal" columns (G, H) plus a small summary block (J, K)
# to the discharge worksheet. Mirrors the existing Q/Qtotal (E/F) layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2: first segment area (measured from bank, i.e. distance 0) ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"

# Running area total (mirrors Qtotal in F2/H-less layout) and the mini
# summary block that echoes the two totals together.
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Row 3: single (non-shared) area formula ---
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- Rows 4-15: shared area formula, same pattern as the D/E columns ---
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Put the cursor where the author left it.
$ws.Range("F5").Select()
